# Update the dSF column (F) values for each data row, per the
# "repull data, push all data, mean calculation" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 6).Value = -3
$ws.Cells.Item(3, 6).Value = -2
$ws.Cells.Item(4, 6).Value = 4
$ws.Cells.Item(5, 6).Value = -6
$ws.Cells.Item(6, 6).Value = 4
$ws.Cells.Item(7, 6).Value = -1
$ws.Cells.Item(8, 6).Value = 2
$ws.Cells.Item(9, 6).Value = 5
$ws.Cells.Item(10, 6).Value = 6
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(12, 6).Value = -4
$ws.Cells.Item(13, 6).Value = -2
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(16, 6).Value = 2
$ws.Cells.Item(17, 6).Value = -6
$ws.Cells.Item(20, 6).Value = 6
$ws.Cells.Item(21, 6).Value = -4
$ws.Cells.Item(23, 6).Value = 6
$ws.Cells.Item(24, 6).Value = -3
$ws.Cells.Item(25, 6).Value = 1
$ws.Cells.Item(26, 6).Value = 1
$ws.Cells.Item(27, 6).Value = 10
$ws.Cells.Item(28, 6).Value = 1
$ws.Cells.Item(30, 6).Value = 1
$ws.Cells.Item(31, 6).Value = 3
$ws.Cells.Item(34, 6).Value = 6
$ws.Cells.Item(37, 6).Value = -1
